# Apply the coordinate_output.xlsx edits:
#  - Rename header cells H1/I1 from "lat"/"lng" to "latitude"/"longitude"
#  - Convert the latitude/longitude data cells (H2:I5) from text strings
#    to real numbers, rounded to 7 decimal places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("H1").Value = "latitude"
$ws.Range("I1").Value = "longitude"

# Data rows 2-5: columns H (latitude) and I (longitude)
for ($row = 2; $row -le 5; $row++) {
    $latCell = $ws.Cells.Item($row, 8)   # column H
    $lngCell = $ws.Cells.Item($row, 9)   # column I

    $lat = [double]$latCell.Text
    $lng = [double]$lngCell.Text

    $latCell.Value = [Math]::Round($lat, 7)
    $lngCell.Value = [Math]::Round($lng, 7)
}
